$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.019.02'
$ws.Range("E2").Value = '  +1.97%  '

$ws.Range("D3").Value = '2.720.47'
$ws.Range("E3").Value = '  +2.25%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '115.16'
$ws.Range("E5").Value = '  +0.97%  '

$ws.Range("D6").Value = '330.31'
$ws.Range("E6").Value = '  +1.14%  '

$ws.Range("D7").Value = '0.530'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  +0.42%  '

$ws.Range("D10").Value = '40.81'
$ws.Range("E10").Value = '  -1.10%  '

$ws.Range("D11").Value = '20.27'
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").Value = '0.0826'
$ws.Range("E12").Value = '  +0.13%  '

$ws.Range("E13").Value = '  +2.38%  '

$ws.Range("D14").Value = '7.62'
$ws.Range("E14").Value = '  +3.24%  '

$ws.Range("D15").Value = '3.141.77'
$ws.Range("E15").Value = '  +2.12%  '

$ws.Range("D16").Value = '2.707.96'
$ws.Range("E16").Value = '  +2.91%  '

$ws.Range("D17").Value = '0.878'
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = '50.896.94'
$ws.Range("E18").Value = '  +1.84%  '

$ws.Range("D19").Value = '13.73'
$ws.Range("E19").Value = '  +3.31%  '

$ws.Range("D20").Value = '2.98'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("D21").Value = '6.83'
$ws.Range("E21").Value = '  +0.41%  '

$ws.Range("D22").Value = '0.0₃0960'
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").Value = '288.10'
$ws.Range("E23").Value = '  +4.10%  '

$ws.Range("D24").Value = '70.14'
$ws.Range("E24").Value = '  -3.34%  '

$ws.Range("D25").Value = '2.60'
$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("D26").Value = '26.82'
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  +2.82%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '35.78'
$ws.Range("E30").Value = '  -2.69%  '

$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.141'
$ws.Range("E31").Value = '  -1.14%  '

$ws.Range("D32").Value = '49.98'
$ws.Range("E32").Value = '  -0.55%  '

$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("D34").Value = '0.0824'
$ws.Range("E34").Value = '  +0.80%  '

$ws.Range("D35").Value = '19.48'
$ws.Range("E35").Value = '  -1.18%  '

$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").Value = '5.02'
$ws.Range("E37").Value = '  -0.53%  '

$ws.Range("D38").Value = '2.09'
$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("E39").Value = '  +2.49%  '

$ws.Range("D40").Value = '23.78'
$ws.Range("E40").Value = '  +6.28%  '

$ws.Range("D41").Value = '128.75'
$ws.Range("E41").Value = '  +3.07%  '

$ws.Range("D42").Value = '0.0352'
$ws.Range("E42").Value = '  +10.20%  '

$ws.Range("D43").Value = '2.31'
$ws.Range("E43").Value = '  +3.88%  '

$ws.Range("E44").Value = '  -0.30%  '

$ws.Range("D45").Value = '3.40'
$ws.Range("E45").Value = '  +1.57%  '

$ws.Range("D46").Value = '2.110.25'
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Value = '2.17'
$ws.Range("E47").Value = '  +8.49%  '

$ws.Range("E48").Value = '  -2.97%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '5.42'
$ws.Range("E49").Value = '  +0.95%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '9.04'
$ws.Range("E50").Value = '  -0.78%  '

$ws.Range("D51").Value = '60.03'
$ws.Range("E51").Value = '  +0.39%  '
